$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete column E ("reviews_count") entirely, shifting columns F:K left to E:J.
$ws.Columns("E").Delete()
